$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.752.88'
$ws.Range("E2").Value = '  +1.93%  '

$ws.Range("D3").Value = '2.797.07'
$ws.Range("E3").Value = '  +1.88%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '351.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '112.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.32%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.557'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.42%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.621'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.66%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.12'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.71%  '

$ws.Range("E11").Value = '  -0.39%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0837'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.85'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.13%  '

$ws.Range("E14").Value = '  +4.31%  '

$ws.Range("D15").Value = '3.238.08'
$ws.Range("E15").Value = '  +1.97%  '

$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.963'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.69%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.801.08'
$ws.Range("E17").Value = '  +1.80%  '

$ws.Range("D18").Value = '51.744.19'
$ws.Range("E18").Value = '  +1.99%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.35'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +10.73%  '

$ws.Range("E20").Value = '  +1.16%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.58'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.14%  '

$ws.Range("E22").Value = '  +2.12%  '

$ws.Range("E23").Value = '  +1.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.55'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.20%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.75'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.47%  '

$ws.Range("E26").Value = '  -0.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.51%  '

$ws.Range("E28").Value = '  +1.08%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '38.74'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +13.79%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.42'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.27'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.61%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '52.65'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.18%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.11'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.25%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0919'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +11.56%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0455'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.19%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.63'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.50%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.90'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.10%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.17'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.95%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.01'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.42%  '

$ws.Range("E41").Value = '  +2.51%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.51'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '122.35'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.67%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.12%  '

$ws.Range("E46").Value = '  +9.41%  '

$ws.Range("E47").Value = '  +8.78%  '

$ws.Range("D48").Value = '2.120.50'
$ws.Range("E48").Value = '  +2.17%  '

$ws.Range("E49").Value = '  +7.95%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.48'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.39%  '

$ws.Range("E51").Value = '  +17.38%  '
